## Generate Report for Archive
## - Status text "Ready for handoff" -> "In Translation" on all sheets
## - Narrow the Status-related columns to match the shorter text's autofit width

$wb = $excel.ActiveWorkbook

# --- 1. Update status text wherever it appears -----------------------------
# Overview sheet: zh-cn / de-de status cells (E2, F2)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# zh-cn sheet: Status column (C2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

# de-de sheet: Status column (C2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- 2. Shrink the columns that held the status text to their new fit ------
# (mirrors the autofit narrowing seen after the shorter text replaced the
# longer "Ready for handoff" string)
$wsOverview.Columns.Item("E").ColumnWidth = 12.5
$wsOverview.Columns.Item("F").ColumnWidth = 12.5
$wsZhCn.Columns.Item("C").ColumnWidth = 12.5
$wsDeDe.Columns.Item("C").ColumnWidth = 12.5
